# Update recomputed NATMI TPM-derived statistics (ligand/receptor expression
# and edge-weight specificity scores) for the Anpep-Sele LR pair sheet, per
# "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.91240233333333
$ws.Range("H2").Value = 35.737207
$ws.Range("I2").Value = 0.04935447906883859
$ws.Range("J2").Value = 0.0493544790688386
$ws.Range("M2").Value = 12.67919733333333
$ws.Range("N2").Value = 38.037592
$ws.Range("O2").Value = 0.9871416146107245
$ws.Range("P2").Value = 0.9871416146107247
$ws.Range("Q2").Value = 151.0396998983938
$ws.Range("R2").Value = 1359.357299085544
$ws.Range("S2").Value = 0.04871986015628454
$ws.Range("T2").Value = 0.04871986015628455
$ws.Range("G3").Value = 11.91240233333333
$ws.Range("H3").Value = 35.737207
$ws.Range("I3").Value = 0.04935447906883859
$ws.Range("J3").Value = 0.0493544790688386
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.01285838538927542
$ws.Range("P3").Value = 0.01285838538927542
$ws.Range("Q3").Value = 1.967424573767889
$ws.Range("R3").Value = 17.706821163911
$ws.Range("S3").Value = 0.0006346189125540536
$ws.Range("T3").Value = 0.0006346189125540537
$ws.Range("I4").Value = 0.5952114870542978
$ws.Range("J4").Value = 0.5952114870542978
$ws.Range("M4").Value = 12.67919733333333
$ws.Range("N4").Value = 38.037592
$ws.Range("O4").Value = 0.9871416146107245
$ws.Range("P4").Value = 0.9871416146107247
$ws.Range("Q4").Value = 1821.527976323414
$ws.Range("R4").Value = 16393.75178691072
$ws.Range("S4").Value = 0.5875580283656299
$ws.Range("T4").Value = 0.58755802836563
$ws.Range("I5").Value = 0.5952114870542978
$ws.Range("J5").Value = 0.5952114870542978
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.01285838538927542
$ws.Range("P5").Value = 0.01285838538927542
$ws.Range("Q5").Value = 23.72699962218667
$ws.Range("S5").Value = 0.007653458688667877
$ws.Range("T5").Value = 0.007653458688667877
$ws.Range("G6").Value = 4.406680666666666
$ws.Range("H6").Value = 13.220042
$ws.Range("I6").Value = 0.01825739449023443
$ws.Range("J6").Value = 0.01825739449023443
$ws.Range("M6").Value = 12.67919733333333
$ws.Range("N6").Value = 38.037592
$ws.Range("O6").Value = 0.9871416146107245
$ws.Range("P6").Value = 0.9871416146107247
$ws.Range("Q6").Value = 55.87317375765156
$ws.Range("R6").Value = 502.858563818864
$ws.Range("S6").Value = 0.01802263387567496
$ws.Range("T6").Value = 0.01802263387567496
$ws.Range("G7").Value = 4.406680666666666
$ws.Range("H7").Value = 13.220042
$ws.Range("I7").Value = 0.01825739449023443
$ws.Range("J7").Value = 0.01825739449023443
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.01285838538927542
$ws.Range("P7").Value = 0.01285838538927542
$ws.Range("Q7").Value = 0.7277970966517778
$ws.Range("R7").Value = 6.550173869866
$ws.Range("S7").Value = 0.0002347606145594679
$ws.Range("T7").Value = 0.0002347606145594679
$ws.Range("G8").Value = 81.382356
$ws.Range("H8").Value = 244.147068
$ws.Range("I8").Value = 0.337176639386629
$ws.Range("J8").Value = 0.3371766393866291
$ws.Range("M8").Value = 12.67919733333333
$ws.Range("N8").Value = 38.037592
$ws.Range("O8").Value = 0.9871416146107245
$ws.Range("P8").Value = 0.9871416146107247
$ws.Range("Q8").Value = 1031.862951175584
$ws.Range("R8").Value = 9286.766560580258
$ws.Range("S8").Value = 0.332841092213135
$ws.Range("T8").Value = 0.3328410922131351
$ws.Range("G9").Value = 81.382356
$ws.Range("H9").Value = 244.147068
$ws.Range("I9").Value = 0.337176639386629
$ws.Range("J9").Value = 0.3371766393866291
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.01285838538927542
$ws.Range("P9").Value = 0.01285838538927542
$ws.Range("Q9").Value = 13.440920024796
$ws.Range("R9").Value = 120.968280223164
$ws.Range("S9").Value = 0.004335547173494016
$ws.Range("T9").Value = 0.004335547173494016
